$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "Intended to give users an easy-to-use music creating
# application " -> "To give users an easy-to-use music creating
# application ".
#
# The visible change is "Intended to" -> "To", which is produced by
# turning "Intended t" into "T" and leaving the rest of the sentence
# untouched. The document's hidden "_GoBack" bookmark (Word's "last
# edit location" marker) ends up sitting right after that new "T",
# since that's where the edit happened.
# ------------------------------------------------------------------
$find1 = $d.Content.Duplicate
$find1.Find.Execute("Intended to give users an easy-to-use music creating application", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find1.Find.Found) { throw "Could not find paragraph 1 target text" }

$introStart = $find1.Start

# "Intended t" is the first 10 characters of the match -> replace with "T"
$replaceRng = $d.Range($introStart, $introStart + 10)
if ($replaceRng.Text -ne "Intended t") { throw "Unexpected text for edit 1: [$($replaceRng.Text)]" }
$replaceRng.Text = "T"

# Relocate the (hidden) "_GoBack" bookmark to sit right after the new "T".
# Adding a bookmark with a name that already exists simply moves it, so
# this both removes it from its old location (after "...bottom right
# corner") and (re)creates it in the new spot.
$bmRng = $d.Range($introStart + 1, $introStart + 1)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ------------------------------------------------------------------
# Edit 2: "he screen displays a 6x4 button board" -> the visible text
# becomes "he screen displays a 4x6 button board", with the "4x6"
# edit splitting the run into "...a 4x6" and " button board".
# ------------------------------------------------------------------
$find2 = $d.Content.Duplicate
$find2.Find.Execute("he screen displays a 6x4 button board", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find2.Find.Found) { throw "Could not find paragraph 2 target text" }

$seg2 = $find2.Duplicate
$seg2.Find.Execute("6x4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $seg2.Find.Found) { throw "Could not find '6x4' in paragraph 2" }

# Temporary bookmarks at the run boundaries (immediately before "he
# screen..." and immediately after "6x4") stop the neighbouring runs
# from being re-merged while the "6x4" text in between is edited.
$tmpBefore = $d.Range($find2.Start, $find2.Start)
$d.Bookmarks.Add("ZZTmpBoundary1", $tmpBefore)

$tmpAfter = $d.Range($seg2.End, $seg2.End)
$d.Bookmarks.Add("ZZTmpBoundary2", $tmpAfter)

$seg2.Text = "4x6"

$d.Bookmarks("ZZTmpBoundary1").Delete()
$d.Bookmarks("ZZTmpBoundary2").Delete()

# ------------------------------------------------------------------
# Edit 3: "Displays the 6x4 button board" -> "Displays the 4x6 button
# board", with the "4x6" edit splitting the run into "...the 4x6" and
# " button board".
# ------------------------------------------------------------------
$find3 = $d.Content.Duplicate
$find3.Find.Execute("Displays the 6x4 button board", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $find3.Find.Found) { throw "Could not find paragraph 3 target text" }

$seg3 = $find3.Duplicate
$seg3.Find.Execute("6x4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $seg3.Find.Found) { throw "Could not find '6x4' in paragraph 3" }

$tmpAfter3 = $d.Range($seg3.End, $seg3.End)
$d.Bookmarks.Add("ZZTmpBoundary3", $tmpAfter3)

$seg3.Text = "4x6"

$d.Bookmarks("ZZTmpBoundary3").Delete()
